# Update the "Förändrad" (Changed) date column (C) for all data rows.
# All rows from 2 to 236 currently hold the serial date 45180 (2023-09-11)
# and must be bumped to 45181 (2023-09-12), keeping formatting/type intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 236
$firstRow = 2

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
